$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new week of data was added at the top of the historical series (weekly
# update). Insert a new row above the current row 74, shifting the existing
# rows 74-104 down to 75-105, and populate the new row with this week's data.
$ws.Rows(74).Insert()

$ws.Cells.Item(74, 1).Value = 10
$ws.Cells.Item(74, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(74, 3).Value = "La Araucanía"
$ws.Cells.Item(74, 4).Value = 44523
$ws.Cells.Item(74, 5).Value = 9
$ws.Cells.Item(74, 6).Value = 100112012
$ws.Cells.Item(74, 7).Value = "Espinaca"
$ws.Cells.Item(74, 8).Value = "Sin especificar"
$ws.Cells.Item(74, 9).Value = "Primera"
$ws.Cells.Item(74, 10).Value = 30
$ws.Cells.Item(74, 11).Value = 8000
$ws.Cells.Item(74, 12).Value = 8000
$ws.Cells.Item(74, 13).Value = 8000
$ws.Cells.Item(74, 14).Value = "`$/docena de atados"
$ws.Cells.Item(74, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(74, 16).Value = 2667
$ws.Cells.Item(74, 17).Value = 3
$ws.Cells.Item(74, 18).Value = "Hortaliza"
